# Translations workbook update: add localization rows for the Arrow and
# Taser objects (object.ARROW / object.GUN) between the existing rows 45
# ("object.JETPACK") and what was row 46 ("Player"-ish / PRIMARY_FIRE
# entries), matching the author's "taser, arrows, other things" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 46 (pushes the former rows 46..126 down to
# 48..128, inheriting the existing formatting of the surrounding rows).
$ws.Rows("46:47").Insert()

# New row 46: object.ARROW / Arrow / <no column C> / Стрілка
$ws.Range("A46").Value2 = "object.ARROW"
$ws.Range("B46").Value2 = "Arrow"
$ws.Range("D46").Value2 = "Стрілка"

# New row 47: object.GUN / Taser / <no column C> / Тазер
$ws.Range("A47").Value2 = "object.GUN"
$ws.Range("B47").Value2 = "Taser"
$ws.Range("D47").Value2 = "Тазер"

# Restore/update the view state (scroll position + active cell) to match
# where the author ended up after editing.
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("D44").Select()
